$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33-101 down to 34-102.
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with its data.
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 44428
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 100114007
$ws.Range("G33").Value = "Jengibre"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = 22083
$ws.Range("N33").Value = "$/caja 13 kilos"
$ws.Range("O33").Value = "Perú"
$ws.Range("P33").Value = 1699
$ws.Range("Q33").Value = 13
$ws.Range("R33").Value = "Hortaliza"

# Match the date-format style used by the other rows in column D.
$ws.Range("D33").NumberFormat = $ws.Range("D34").NumberFormat
